$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Classi" (sheet1.xml): add week row 16 for "18-23 aprile"
# ---------------------------------------------------------------------------
$wsClassi = $wb.Worksheets.Item("Classi")
$wsClassi.Range("A16").Value = "18-23 aprile"
$wsClassi.Range("B16").Value = 4625
$wsClassi.Range("B16").NumberFormat = "#,##0"
$wsClassi.Range("C16").Value = 8157
$wsClassi.Range("C16").NumberFormat = "#,##0"
$wsClassi.Range("D16").Value = 0.56700000000000006
$wsClassi.Range("D16").NumberFormat = "0.0%"
$wsClassi.Range("E16").Value = 376584
$wsClassi.Range("E16").NumberFormat = "#,##0"
$wsClassi.Range("F16").Value = 213765
$wsClassi.Range("F16").NumberFormat = "#,##0"
$wsClassi.Range("G16").Value = 0.56799999999999995
$wsClassi.Range("G16").NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Sheet "Alunni in presenza" (sheet2.xml): add week row 16
# ---------------------------------------------------------------------------
$wsAlunniPresenza = $wb.Worksheets.Item("Alunni in presenza")
$wsAlunniPresenza.Range("A16").Value = "18-23 aprile"
$wsAlunniPresenza.Range("B16").Value = 7396217
$wsAlunniPresenza.Range("B16").NumberFormat = "#,##0"
$wsAlunniPresenza.Range("C16").Value = 4184506
$wsAlunniPresenza.Range("C16").NumberFormat = "#,##0"
$wsAlunniPresenza.Range("D16").Value = 0.56600000000000006
$wsAlunniPresenza.Range("D16").NumberFormat = "0.0%"
$wsAlunniPresenza.Range("E16").Value = 4132398
$wsAlunniPresenza.Range("E16").NumberFormat = "#,##0"
$wsAlunniPresenza.Range("F16").Value = 0.98799999999999999
$wsAlunniPresenza.Range("F16").NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Sheet "Alunni" (sheet3.xml): add rows 44-46 (Infanzia / Primaria / Sec.)
# for "18-23 aprile" and turn E38:E46 into one shared formula
# ---------------------------------------------------------------------------
$wsAlunni = $wb.Worksheets.Item("Alunni")

$wsAlunni.Range("A44").Value = "18-23 aprile"
$wsAlunni.Range("B44").Value = "Infanzia"
$wsAlunni.Range("C44").Value = 484017
$wsAlunni.Range("C44").NumberFormat = "#,##0"
$wsAlunni.Range("D44").Value = 480691
$wsAlunni.Range("D44").NumberFormat = "#,##0"
$wsAlunni.Range("F44").Value = 0.99299999999999999
$wsAlunni.Range("F44").NumberFormat = "0.0%"

$wsAlunni.Range("A45").Value = "18-23 aprile"
$wsAlunni.Range("B45").Value = "Primaria"
$wsAlunni.Range("C45").Value = 1342673
$wsAlunni.Range("C45").NumberFormat = "#,##0"
$wsAlunni.Range("D45").Value = 1325560
$wsAlunni.Range("D45").NumberFormat = "#,##0"
$wsAlunni.Range("F45").Value = 0.98699999999999999
$wsAlunni.Range("F45").NumberFormat = "0.0%"

$wsAlunni.Range("A46").Value = "18-23 aprile"
$wsAlunni.Range("B46").Value = "Sec. 1° e 2° Grado"
$wsAlunni.Range("C46").Value = 2357816
$wsAlunni.Range("C46").NumberFormat = "#,##0"
$wsAlunni.Range("D46").Value = 2326147
$wsAlunni.Range("D46").NumberFormat = "#,##0"
$wsAlunni.Range("F46").Value = 0.98699999999999999
$wsAlunni.Range("F46").NumberFormat = "0.0%"

# Re-apply the formula across the whole column so it becomes a single
# shared formula group E38:E46 (matches the original authoring pattern).
$wsAlunni.Range("E38:E46").Formula = "=C38-D38"
$wsAlunni.Range("E38:E46").NumberFormat = "#,##0"

# ---------------------------------------------------------------------------
# Sheet "Personale scolastico" (sheet4.xml): add week row 16
# ---------------------------------------------------------------------------
$wsPersonale = $wb.Worksheets.Item("Personale scolastico")
$wsPersonale.Range("A16").Value = "18-23 aprile"
$wsPersonale.Range("B16").Value = 775867
$wsPersonale.Range("B16").NumberFormat = "#,##0"
$wsPersonale.Range("C16").Value = 436867
$wsPersonale.Range("C16").NumberFormat = "#,##0"
$wsPersonale.Range("D16").Value = 0.56299999999999994
$wsPersonale.Range("D16").NumberFormat = "0.0%"
$wsPersonale.Range("E16").Value = 424494
$wsPersonale.Range("E16").NumberFormat = "#,##0"
$wsPersonale.Range("F16").Value = 0.97199999999999998
$wsPersonale.Range("F16").NumberFormat = "0.0%"
$wsPersonale.Range("G16").Value = 204526
$wsPersonale.Range("G16").NumberFormat = "#,##0"
$wsPersonale.Range("H16").Value = 115592
$wsPersonale.Range("H16").NumberFormat = "#,##0"
$wsPersonale.Range("I16").Value = 0.56499999999999995
$wsPersonale.Range("I16").NumberFormat = "0.0%"
$wsPersonale.Range("J16").Value = 113013
$wsPersonale.Range("J16").NumberFormat = "#,##0"
$wsPersonale.Range("K16").Value = 0.97799999999999998
$wsPersonale.Range("K16").NumberFormat = "0.0%"

# ---------------------------------------------------------------------------
# Restore the selections shown in each sheet's view after the edits.
# "Personale scolastico" is activated last so it remains the active tab,
# matching the workbook's activeTab setting.
# ---------------------------------------------------------------------------
$wsClassi.Activate()
$wsClassi.Range("A16").Select()

$wsAlunniPresenza.Activate()
$wsAlunniPresenza.Range("A16").Select()

$wsAlunni.Activate()
$wsAlunni.Range("C47").Select()

$wsPersonale.Activate()
$wsPersonale.Range("E17").Select()
